$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "[What : introduce your research questions and hypothesis.]"
#    Collapse the split runs "introduce your research " / "questions
#    and hypothesis." / "]" into a single run by replacing the text
#    with itself (Word re-merges the run when old == new).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "introduce your research questions and hypothesis.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "introduce your research questions and hypothesis.", 2)

# ---------------------------------------------------------------------
# 2) "[What : introduce the research methods and data sources you used
#    for the analysis]" -- same trick.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "introduce the research methods and data sources you used for the analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "introduce the research methods and data sources you used for the analysis", 2)

# ---------------------------------------------------------------------
# 3) Merge the two places where a curly close-quote run "<E2><80><9D>"
#    is immediately followed by a separate ")" run: "Speaking”)" and
#    "…Children”) ". Looping the same find/replace lets Word walk
#    forward through every match in the document; the two occurrences
#    that are already a single run (”) / ”). / ”).) are no-ops.
# ---------------------------------------------------------------------
$rngQuote = $d.Content
$guard = 0
while ($rngQuote.Find.Execute("”)", $true, $false, $false, $false, $false,
                               $true, 1, $false, "”)", 2)) {
    $rngQuote.Collapse(0)
    $guard++
    if ($guard -gt 20) { break }
}

# ---------------------------------------------------------------------
# 4) "-Results " -> "Results " (drop the leading hyphen). The text is
#    unique in the document (the two "-Discussion" bullets are left
#    untouched).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("-Results ", $true, $false, $false, $false,
                                 $false, $true, 1, $false, "Results ", 2)

# ---------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the References paragraph to right
#    after the ")" that closes "-Discussion (power of
#    narrative//storytelling)" (the second such bullet, numbered list
#    id 2). A collapsed Range sitting exactly on that paragraph-mark
#    boundary trips up Bookmarks.Add, so a small sentinel is inserted,
#    used to anchor the bookmark, then removed again.
# ---------------------------------------------------------------------
$rngTarget = $d.Content
$null = $rngTarget.Find.Execute(
    "-Discussion (power of narrative//storytelling)", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
$rngTarget.Collapse(0)
$null = $rngTarget.Find.Execute(
    "-Discussion (power of narrative//storytelling)", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
$rngTarget.Collapse(0)
$rngTarget.InsertAfter("ZZZGoBackSentinelZZZ")

$rngSentinel = $d.Content
$null = $rngSentinel.Find.Execute("ZZZGoBackSentinelZZZ", $true, $false,
                                   $false, $false, $false, $true, 1,
                                   $false, "", 0)
$rngSentinel.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rngSentinel)

$rngCleanup = $d.Content
$null = $rngCleanup.Find.Execute("ZZZGoBackSentinelZZZ", $true, $false,
                                  $false, $false, $false, $true, 1,
                                  $false, "", 0)
$rngCleanup.Delete()

Write-Output "done"
